$d = $word.ActiveDocument

# Find the "Requisitos" Heading2 paragraph and remove it, together with
# the "LOQ4205 - Sistemas Produtivos II (Requisito fraco)" bullet
# paragraph that follows it, all the way to the end of the document body
# (just before the final section break).
$startRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.TrimEnd("`r", "`a") -eq "Requisitos") {
        $startRange = $p.Range.Start
        break
    }
}

if ($startRange -ne $null) {
    $endRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
    $rng = $d.Range($startRange, $endRange)
    $rng.Delete()
}
